$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
# Force B2 to be stored as text "1" (not a number) while keeping its original (no) style
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = 15
$ws.Range("E2").Value = "2024-04-22 11:12:56"

# Delete rows 3 through 5 (they are removed entirely in the target)
$ws.Range("A3:E5").EntireRow.Delete()
